# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# OFF sheet (Target Depth Data - offense): row 3 ("R") updated with new cumulative totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 528
$wsOff.Range("C3").Value = 367
$wsOff.Range("D3").Value = 111
$wsOff.Range("E3").Value = 59
$wsOff.Range("F3").Value = 13

# DEF sheet (Target Depth Data - defense): row 3 ("R") updated with new cumulative totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 428
$wsDef.Range("C3").Value = 298
$wsDef.Range("D3").Value = 115
$wsDef.Range("E3").Value = 67
